$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.664.02"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "3.340.16"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.34"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.72"
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.647"
$ws.Range("E7").Value = "  +3.41%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "3.337.77"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("E11").Value = "  +2.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.406"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").Value = "3.917.99"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").Value = "66.597.27"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.67"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000166"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.316.84"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "427.93"
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.21"
$ws.Range("E21").Value = "  -2.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.41"
$ws.Range("E22").Value = "  -2.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.21"
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.68"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").Value = "3.467.41"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.519"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  +5.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000116"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.06"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.50"
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.26"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.66"
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.20"
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.62"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.45"
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.82"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("D41").Value = "2.873.31"
$ws.Range("E41").Value = "  +2.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.58"
$ws.Range("E42").Value = "  -5.18%  "
$ws.Range("E43").Value = "  -2.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.763"
$ws.Range("E44").Value = "  -4.72%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0668"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.85"
$ws.Range("E46").Value = "  -1.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.02"
$ws.Range("E47").Value = "  -3.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.33"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.35"
$ws.Range("E49").Value = "  -4.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "315.21"
$ws.Range("E50").Value = "  -3.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0274"
$ws.Range("E51").Value = "  +0.34%  "
